$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the shortened descriptions in column C (Week1)
$ws.Range("C4").Value = "Mô tả nghiệp vụ Admin, vẽ sơ đồ usecase Admin"
$ws.Range("C5").Value = "Mô tả nghiệp vụ User, vẽ sơ đồ usecase User"

# Add new header for week3 in column E
$ws.Range("E3").Value = "week3"

# Add new content for column D (new Week "in-between" column)
$ws.Range("D4").Value = "Mô ta usecae Admin,  vẽ sơ đồ use case, ghi vào file Usecase specification"
$ws.Range("D5").Value = "mô tả use case học viên,  vẽ sơ đồ use case, ghi vào file Usecase specification"

# Add new content for column E rows 4 and 5 (same text for both)
$ws.Range("E4").Value = "Vẽ flow use case, và bảng Flow UC"
$ws.Range("E5").Value = "Vẽ flow use case, và bảng Flow UC"

# Set new column D width to match target (~15.42578125 chars stored);
# the COM width setter here quantizes to 1/6-character increments, so
# 14.6 is the closest input that lands on the nearest representable width (15.5).
$ws.Columns.Item(4).ColumnWidth = 14.6
